$wb = $excel.ActiveWorkbook

# NOTE: the workbook contains two worksheets whose names differ only by
# case ("Vector_bf" and "Vector_BF"). Worksheets.Item(<name>) resolves
# case-insensitively, so those two sheets are addressed by their
# 1-based position in the workbook instead, to make sure each edit
# lands on the intended sheet.
$sheet5Name = $wb.Worksheets.Item(5).Name
$sheet6Name = $wb.Worksheets.Item(6).Name
if ($sheet5Name -ne "Vector_bf") { throw "Expected worksheet #5 to be 'Vector_bf' but found '$sheet5Name'" }
if ($sheet6Name -ne "Vector_BF") { throw "Expected worksheet #6 to be 'Vector_BF' but found '$sheet6Name'" }

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("A2").Value = "-12.266066838046278 + 1.053984575835476y_1 + 1.208226221079692y_2"
$ws.Range("B2").Value = "12.266066838046278"
$ws.Range("D2").Value = "0.75"
$ws.Range("E2").Value = "4.3"
$ws.Range("F2").Value = "3.5999999999999996"
$ws.Range("A3").Value = "-16.59125964010283 + 1.7866323907455013y_1 + 0.46272493573264795y_2"
$ws.Range("B3").Value = "12.59125964010283"
$ws.Range("D3").Value = "0.19"
$ws.Range("E3").Value = "1.9"
$ws.Range("F3").Value = "4.2"
$ws.Range("A4").Value = "54.89266508652872 - 2x - 4.895630031679144y_1 + 0.5319823343063854y_2"
$ws.Range("B4").Value = "-70.89266508652872"
$ws.Range("D4").Value = "0.24"
$ws.Range("E4").Value = "9.8"
$ws.Range("F4").Value = "1.2"
$ws.Range("A5").Value = "-65.30201371036847 + 8x + 0.9854327335047129y_1 - 0.00856898029134534y_2"
$ws.Range("B5").Value = "17.252013710368466"
$ws.Range("D5").Value = "0.92"
$ws.Range("E5").Value = "2.4"
$ws.Range("F5").Value = "2.0"
$ws.Range("A6").Value = "11.024293059125966 - 2x - 0.1863753213367607y_1 + 1.0668380462724938y_2"
$ws.Range("B6").Value = "-0.9757069408740331"
$ws.Range("D6").Value = "0.82"
$ws.Range("E6").Value = "4.8"
$ws.Range("F6").Value = "0.8999999999999999"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("A2").Value = "7.1"
$ws.Range("B2").Value = "8.600000000000001"
$ws.Range("C2").Value = "2.65"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A2").Value = "0.2912322701565502"
$ws.Range("A3").Value = "-1.9886868999079117"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A2").Value = "11.000000000000002"
$ws.Range("A3").Value = "41.58000207395175"
$ws.Range("A4").Value = "-18.388234074146023"

# --- Vector_Alpha ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 2.04
$ws.Range("A3").Value = 1.2000000000000002

